$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 values: Day (serial date) through hourly prices, daily avg, and slot summary columns
$ws.Range("A2").Value = 46034
$ws.Range("B2").Value = 77.48999999999999
$ws.Range("C2").Value = 77.94
$ws.Range("D2").Value = 73.01000000000001
$ws.Range("E2").Value = 63.43
$ws.Range("F2").Value = 60.13
$ws.Range("G2").Value = 67.89
$ws.Range("H2").Value = 82.45
$ws.Range("I2").Value = 97.14
$ws.Range("J2").Value = 104.6
$ws.Range("K2").Value = 101.29
$ws.Range("L2").Value = 89.33
$ws.Range("M2").Value = 78.75
$ws.Range("N2").Value = 76.37
$ws.Range("O2").Value = 76.87
$ws.Range("P2").Value = 77.33
$ws.Range("Q2").Value = 76.03
$ws.Range("R2").Value = 85.62
$ws.Range("S2").Value = 102.85
$ws.Range("T2").Value = 116.49
$ws.Range("U2").Value = 115.7
$ws.Range("V2").Value = 119.64
$ws.Range("W2").Value = 118.42
$ws.Range("X2").Value = 99.23999999999999
$ws.Range("Y2").Value = 86.31
$ws.Range("Z2").Value = 88.51000000000001

$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 105.9
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 119.03
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 116.1
$ws.Range("AG2").Value = "0h-23h"
